$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns (title_en / description_en) right after
#     description_vi (old column F), shifting everything from the old
#     "gen_map_type" column onward two slots to the right. ---
$ws.Columns("G:H").Insert()

# --- Header row ---
$ws.Range("G1").Value = "title_en"
$ws.Range("H1").Value = "description_en"

# --- New English title/description content for each requirement row ---
$ws.Range("G2").Value = "Wall-Following Algorithm"
$ws.Range("H2").Value = "This maze is quite complex! Apply the 'wall-following' algorithm (always keep one wall on your right) to find the exit."

$ws.Range("G3").Value = "Flood-Fill Algorithm"
$ws.Range("H3").Value = "Your task is to collect all the gems on the islands. Design an algorithm to explore every corner without missing any spot."

$ws.Range("G4").Value = "Systematic Collection"
$ws.Range("H4").Value = "Gems are scattered throughout a large maze. You need an efficient algorithm to ensure none are left behind."

$ws.Range("G5").Value = "Shortest Path"
$ws.Range("H5").Value = "There are many ways to reach the goal, but only one is the most efficient. Write an algorithm to find the optimal route."

$ws.Range("G6").Value = "The Delivery Challenge"
$ws.Range("H6").Value = "A real test awaits! Collect all different types of treasures before reaching the goal. The order of collection is up to you!"

# --- Match the source/description columns' wrap-text formatting ---
$ws.Range("F1:F6").WrapText = $true
$ws.Range("G2:G6").WrapText = $true
$ws.Range("H1:H6").WrapText = $true

# --- Leftover formatted-but-empty cells (e.g. from a paste/fill action) ---
$ws.Range("G16:J20").Font.Name = "Helvetica"

# --- Final selection left on the new description_en header cell ---
$ws.Range("H1").Select() | Out-Null
